$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FluxData")

# Insert a new row at position 8 (shifts old rows 8..41 down to 9..42),
# so the existing DIL_* entries move down by one row and a new entry
# "EX_o2_e.f" can be placed at row 8, matching the other EX_* rows above it.
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the EX_o2_e.f flux entry.
$ws.Cells.Item(8, 1).Value = "EX_o2_e.f"
$ws.Cells.Item(8, 2).Value = 20.36
$ws.Cells.Item(8, 3).Formula = "=MIN(0.92, 0.1*B8)"

# Match formatting of the other EX_* rows (rows 2-7): row height 13.8,
# default (non-wrapping) style instead of the wrapping style used by DIL_* rows.
$ws.Rows.Item(8).RowHeight = 13.8

# The row that used to be the last DIL_* row (old row 41, "DIL_amp_d1.f") has now
# shifted down to row 42; re-create that trailing row with its original content.
$ws.Cells.Item(42, 1).Value = "DIL_amp_d1.f"
$ws.Cells.Item(42, 2).Value = 100
$ws.Cells.Item(42, 3).Value = 0.0001

# Restore the selection to A1 as left by the edit.
$ws.Range("A1").Select()
